$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 1).Value = "ar2"
$ws.Cells.Item(1, 2).Value = "A/C"
$ws.Cells.Item(1, 3).Value = 30
$ws.Cells.Item(1, 4).Value = $true

# Row 2
$ws.Cells.Item(2, 1).Value = "a"
$ws.Cells.Item(2, 2).Value = "Lâmpada"
$ws.Cells.Item(2, 3).Value = 100
$ws.Cells.Item(2, 4).Value = $false

# Row 3
$ws.Cells.Item(3, 1).Value = "ae"
$ws.Cells.Item(3, 2).Value = "A/C"
$ws.Cells.Item(3, 3).Value = 23
$ws.Cells.Item(3, 4).Value = $false

# Row 4
$ws.Cells.Item(4, 1).Value = "la"
$ws.Cells.Item(4, 2).Value = "Lâmpada"
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = $false

# Row 5
$ws.Cells.Item(5, 1).Value = "tv"
$ws.Cells.Item(5, 2).Value = "Televisor"
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = $false

# Row 6
# A6 holds the digit-string "6" (must stay text, not be coerced to a number)
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "6"
$ws.Cells.Item(6, 1).Style = "Normal"
$ws.Cells.Item(6, 2).Value = "A/C"
$ws.Cells.Item(6, 3).Value = 23
$ws.Cells.Item(6, 4).Value = $false
